$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row for LARISSA (004289402 / 349.86) above existing row 87 (004907688)
$ws.Rows.Item(87).Insert()
$ws.Cells.Item(87, 1).NumberFormat = "@"
$ws.Cells.Item(87, 1).Value = "004289402"
$ws.Cells.Item(87, 2).Value = "LARISSA"
$ws.Cells.Item(87, 3).Value = 349.86

# Insert new row for LARISSA (004290978 / 38.7) above existing row 194 (003435941)
# (the earlier insert shifted everything below row 87 down by one, so the
# original row 194 is now row 195)
$ws.Rows.Item(195).Insert()
$ws.Cells.Item(195, 1).NumberFormat = "@"
$ws.Cells.Item(195, 1).Value = "004290978"
$ws.Cells.Item(195, 2).Value = "LARISSA"
$ws.Cells.Item(195, 3).Value = 38.7

# Remove the old LARISSA row (004290978 / -111.3), originally row 289, now
# shifted down by the two inserts above to row 291
$ws.Rows.Item(291).Delete()

# Remove the MARCELO row (004748761 / -28536.57), originally row 291, now
# shifted down by the two inserts above (and one fewer after the prior
# delete) to row 292
$ws.Rows.Item(292).Delete()
